$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @{Row=2; D="27.821.21"; E="  +0.48%  "},
  @{Row=3; D="1.907.18"; E="  +0.68%  "},
  @{Row=4; D="0.9994"; E="  -0.18%  "},
  @{Row=5; D="312.82"; E="  +0.02%  "},
  @{Row=6; D="0.9996"; E="  -0.09%  "},
  @{Row=7; E="  +5.95%  "},
  @{Row=8; D="0.3787"; E="  -0.38%  "},
  @{Row=9; D="0.07236"; E="  -1.21%  "},
  @{Row=10; D="21.32"; E="  +3.68%  "},
  @{Row=11; D="0.9082"; E="  -0.65%  "},
  @{Row=12; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.933.30"; E="  +1.84%  "},
  @{Row=13; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.07638"; E="  -0.46%  "},
  @{Row=14; D="5.452"; E="  -0.41%  "},
  @{Row=15; D="92.27"; E="  +1.30%  "},
  @{Row=16; D="1.0000"; E="  -0.18%  "},
  @{Row=17; D="0.000008705"; E="  -0.80%  "},
  @{Row=18; D="0.9992"; E="  -0.10%  "},
  @{Row=19; D="27.840.68"; E="  -0.07%  "},
  @{Row=20; D="14.54"; E="  +0.14%  "},
  @{Row=21; D="5.150"; E="  +0.47%  "},
  @{Row=22; D="2.165.81"; E="  +0.06%  "},
  @{Row=23; D="10.85"; E="  +0.98%  "},
  @{Row=24; D="6.624"; E="  +0.37%  "},
  @{Row=25; D="154.00"; E="  +0.39%  "},
  @{Row=26; D="1.868"; E="  -2.19%  "},
  @{Row=27; D="2.171"; E="  +1.31%  "},
  @{Row=28; D="18.32"; E="  -0.30%  "},
  @{Row=29; D="114.50"; E="  -0.99%  "},
  @{Row=30; D="4.851"; E="  -0.89%  "},
  @{Row=31; D="0.09023"; E="  +0.95%  "},
  @{Row=32; D="4.865"; E="  +4.80%  "},
  @{Row=33; D="3.178"; E="  -0.50%  "},
  @{Row=34; D="1.231"; E="  +0.94%  "},
  @{Row=35; D="0.7816"; E="  +2.06%  "},
  @{Row=36; D="0.02093"; E="  +3.12%  "},
  @{Row=37; D="2.620"; E="  +3.90%  "},
  @{Row=38; D="3.076"; E="  +3.16%  "},
  @{Row=39; D="1.094"; E="  +0.05%  "},
  @{Row=40; D="0.5550"; E="  +1.49%  "},
  @{Row=41; D="0.05282"; E="  +0.07%  "},
  @{Row=42; D="6.693"; E="  -3.08%  "},
  @{Row=43; D="114.54"; E="  +1.78%  "},
  @{Row=44; D="8.579"; E="  +0.62%  "},
  @{Row=45; D="0.1513"; E="  -0.24%  "},
  @{Row=46; D="0.4810"; E="  +0.40%  "},
  @{Row=47; E="  -1.41%  "},
  @{Row=48; D="0.9994"; E="  -0.07%  "},
  @{Row=49; D="1.621"; E="  -0.61%  "},
  @{Row=50; D="66.92"; E="  -0.77%  "},
  @{Row=51; D="0.05997"; E="  -0.86%  "}
)

foreach ($item in $updates) {
  if ($item.ContainsKey("B")) { $ws.Cells.Item($item.Row, 2).NumberFormat = "@"; $ws.Cells.Item($item.Row, 2).Value = $item.B }
  if ($item.ContainsKey("C")) { $ws.Cells.Item($item.Row, 3).NumberFormat = "@"; $ws.Cells.Item($item.Row, 3).Value = $item.C }
  if ($item.ContainsKey("D")) { $ws.Cells.Item($item.Row, 4).NumberFormat = "@"; $ws.Cells.Item($item.Row, 4).Value = $item.D }
  if ($item.ContainsKey("E")) { $ws.Cells.Item($item.Row, 5).NumberFormat = "@"; $ws.Cells.Item($item.Row, 5).Value = $item.E }
}
